# Apply updates to the kit_equipamento sheet ("criando a separacao de miscelaneas")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 2
$ws.Range("C3").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("C11").Value = 5
$ws.Range("C14").Value = 3
$ws.Range("C15").Value = 1
$ws.Range("C21").Value = 1
$ws.Range("C23").Value = 3
